$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TotalEnergy")
$ws2 = $wb.Worksheets.Item("Meters")

# Update TotalEnergy sheet (sheet1): rewrite rows 2-6
$ws.Range("A2").Value = "Total Site Energy [GJ]"
$ws.Range("A3").Value = "Natural Gas Total End Uses [GJ]"
$ws.Range("A4").Value = "Electricity Total End Uses [GJ]"
$ws.Range("A5").Value = "Total Source Energy [GJ]"
$ws.Range("A6").Value = "Electricity Heat Rejection [GJ]"

# Update selections to match the diff
$ws2.Range("A3").Select()
$ws.Range("A6").Select()
